$d = $word.ActiveDocument

# --- Create the three new character styles ---
$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Datas da campaña de 2022 ..." run (4 occurrences) ---
$datesText = "Datas da campaña de 2022 que usan constelación de Orión: 16-25 de xaneiro, 14-23 de febreiro, 14-24 de marzo"
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($found) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $found = $rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

# --- Apply GaNParagraph to the "Grazas por participar ..." run ---
$thanksText = "Grazas por participar nesta campaña global de medida da contaminación lumínica mediante a observación das estrelas máis febles que podes albiscar. Localizando e observando a  constelación de Orión e comparándoa co que aparece nos mapas estelares recollidos neste documento podes saber canto contribúen á contaminación lumínica os sistemas de iluminación que hai no teu barrio ou vila. As túas achegas á base de datos en liña de GLOBE at Night (O MUNDO á Noite) servirán para documentar a calidade do ceo nocturno."
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute($thanksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Os mapas de estrelas ..." run ---
$mapsText = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$found3 = $rng3.Find.Execute($mapsText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Style = "GaNLinks"
}

Write-Output "Applied GaNStyle/GaNParagraph/GaNLinks styles."
